$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1218.237060546875
$ws.Range("C2").Value = 0.9139
$ws.Range("D2").Value = 0.9157000184059143
$ws.Range("E2").Value = 1.326300024986267
$ws.Range("F2").Value = 0.5130000114440918
$ws.Range("H2").Value = 0.7594

$ws.Range("B3").Value = 1164.693115234375
$ws.Range("C3").Value = 0.9318
$ws.Range("D3").Value = 0.9326
$ws.Range("E3").Value = 1.04040002822876
$ws.Range("F3").Value = 0.7684999704360962
$ws.Range("H3").Value = 0.9093

$ws.Range("B4").Value = 797.730712890625
$ws.Range("C4").Value = 0.9319
$ws.Range("D4").Value = 0.9347
$ws.Range("E4").Value = 1.037099957466125
$ws.Range("F4").Value = 0.7386999726295471
$ws.Range("H4").Value = 0.9274

$ws.Range("B5").Value = 828.351806640625
$ws.Range("C5").Value = 0.8738
$ws.Range("D5").Value = 0.8801
$ws.Range("E5").Value = 0.9781000018119812
$ws.Range("F5").Value = 0.5246999859809875
$ws.Range("H5").Value = 0.4439

$ws.Range("B6").Value = 1137.451171875
$ws.Range("C6").Value = 0.8999
$ws.Range("D6").Value = 0.9089
$ws.Range("E6").Value = 0.9739000201225281
$ws.Range("F6").Value = 0.6384000182151794
$ws.Range("H6").Value = 0.6989

$ws.Range("B7").Value = 890.1533813476562
$ws.Range("C7").Value = 0.8964
$ws.Range("D7").Value = 0.9004999995231628
$ws.Range("E7").Value = 0.9761999845504761
$ws.Range("F7").Value = 0.7508000135421753
$ws.Range("H7").Value = 0.6245000000000001

$ws.Range("B8").Value = 985.451171875
$ws.Range("C8").Value = 0.883
$ws.Range("D8").Value = 0.8844
$ws.Range("E8").Value = 0.9664000272750854
$ws.Range("F8").Value = 0.7779999971389771
$ws.Range("H8").Value = 0.4823

$ws.Range("B9").Value = 7022.068359375
$ws.Range("C9").Value = 0.9049
$ws.Range("D9").Value = 0.9121
$ws.Range("E9").Value = 1.326300024986267
$ws.Range("F9").Value = 0.5130000114440918
$ws.Range("H9").Value = 4.845700000000001
